$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "Before delete A37:" $ws.Range("A37").Value2 "B37:" $ws.Range("B37").Value2
$ws.Rows(37).Delete()
Write-Host "After delete A37:" $ws.Range("A37").Value2 "B37:" $ws.Range("B37").Value2
